$d = $word.ActiveDocument

# 1) Replace the changed portion of the phone number. Word, when a user
#    selects "19)523-5637" and types "84)218-4109" over it, keeps the
#    untouched prefix "Tel: (9" as the original run and puts the newly
#    typed text in a fresh run.
$find = $d.Content.Find
$find.Execute("19)523-5637", $true, $false, $false, $false, $false, `
               $true, 1, $false, "84)218-4109", 2)

# 2) Word also re-anchors the hidden "_GoBack" bookmark (tracks the most
#    recent edit position) to the spot of this edit, removing it from
#    its previous location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$editedRange = $d.Content
$editedRange.Find.Execute("84)218-4109")
$goBackRange = $d.Range($editedRange.End, $editedRange.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)
